# Updates the cryptos price list (Coin/Link/Price/Volume(1h)) to the latest
# scraped snapshot. Price values (column D) are written with a leading
# apostrophe so Excel stores them as literal text (preserving exact
# formatting such as trailing zeros / grouped price strings) instead of
# reinterpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''29.179.99'
$ws.Cells.Item(2, 5).Value = '  +0.98%  '

$ws.Cells.Item(3, 4).Value = '''1.936.97'
$ws.Cells.Item(3, 5).Value = '  +2.28%  '

$ws.Cells.Item(4, 4).Value = '''0.9990'
$ws.Cells.Item(4, 5).Value = '  -0.46%  '

$ws.Cells.Item(5, 4).Value = '''326.16'
$ws.Cells.Item(5, 5).Value = '  -0.06%  '

$ws.Cells.Item(6, 4).Value = '''0.9994'
$ws.Cells.Item(6, 5).Value = '  -0.34%  '

$ws.Cells.Item(7, 4).Value = '''0.4616'
$ws.Cells.Item(7, 5).Value = '  +0.55%  '

$ws.Cells.Item(8, 4).Value = '''0.3902'
$ws.Cells.Item(8, 5).Value = '  -0.03%  '

$ws.Cells.Item(9, 4).Value = '''0.07871'
$ws.Cells.Item(9, 5).Value = '  +0.39%  '

$ws.Cells.Item(10, 4).Value = '''0.9967'
$ws.Cells.Item(10, 5).Value = '  +0.80%  '

$ws.Cells.Item(11, 4).Value = '''22.17'
$ws.Cells.Item(11, 5).Value = '  +1.21%  '

$ws.Cells.Item(12, 4).Value = '''1.916.94'
$ws.Cells.Item(12, 5).Value = '  +0.60%  '

$ws.Cells.Item(13, 4).Value = '''5.830'
$ws.Cells.Item(13, 5).Value = '  +2.28%  '

$ws.Cells.Item(14, 4).Value = '''7.095'
$ws.Cells.Item(14, 5).Value = '  +0.78%  '

$ws.Cells.Item(15, 5).Value = '  +1.38%  '

$ws.Cells.Item(16, 4).Value = '''87.74'
$ws.Cells.Item(16, 5).Value = '  -0.37%  '

$ws.Cells.Item(17, 4).Value = '''1.002'
$ws.Cells.Item(17, 5).Value = '  -0.23%  '

$ws.Cells.Item(18, 4).Value = '''0.000009960'
$ws.Cells.Item(18, 5).Value = '  +0.05%  '

$ws.Cells.Item(19, 4).Value = '''17.12'
$ws.Cells.Item(19, 5).Value = '  +1.15%  '

$ws.Cells.Item(20, 4).Value = '''1.001'

$ws.Cells.Item(21, 4).Value = '''29.238.04'
$ws.Cells.Item(21, 5).Value = '  +1.18%  '

$ws.Cells.Item(22, 4).Value = '''5.507'
$ws.Cells.Item(22, 5).Value = '  +4.03%  '

$ws.Cells.Item(23, 5).Value = '  +2.23%  '

$ws.Cells.Item(24, 4).Value = '''2.175.77'
$ws.Cells.Item(24, 5).Value = '  +3.71%  '

$ws.Cells.Item(25, 4).Value = '''2.098'
$ws.Cells.Item(25, 5).Value = '  +1.67%  '

$ws.Cells.Item(26, 4).Value = '''155.73'
$ws.Cells.Item(26, 5).Value = '  -0.18%  '

$ws.Cells.Item(27, 4).Value = '''19.47'
$ws.Cells.Item(27, 5).Value = '  +0.82%  '

$ws.Cells.Item(28, 4).Value = '''5.910'
$ws.Cells.Item(28, 5).Value = '  -0.33%  '

$ws.Cells.Item(29, 4).Value = '''118.76'
$ws.Cells.Item(29, 5).Value = '  +0.93%  '

$ws.Cells.Item(30, 4).Value = '''1.879'
$ws.Cells.Item(30, 5).Value = '  -2.70%  '

$ws.Cells.Item(31, 4).Value = '''0.09334'
$ws.Cells.Item(31, 5).Value = '  -0.16%  '

$ws.Cells.Item(32, 4).Value = '''0.8922'
$ws.Cells.Item(32, 5).Value = '  -2.28%  '

$ws.Cells.Item(33, 4).Value = '''5.218'
$ws.Cells.Item(33, 5).Value = '  -1.39%  '

$ws.Cells.Item(34, 4).Value = '''1.330'
$ws.Cells.Item(34, 5).Value = '  -0.38%  '

$ws.Cells.Item(35, 4).Value = '''3.138'
$ws.Cells.Item(35, 5).Value = '  -4.39%  '

$ws.Cells.Item(36, 4).Value = '''0.05787'
$ws.Cells.Item(36, 5).Value = '  +0.40%  '

$ws.Cells.Item(37, 4).Value = '''1.171'
$ws.Cells.Item(37, 5).Value = '  -1.67%  '

$ws.Cells.Item(38, 2).Value = 'PEPE'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(38, 4).Value = '''0.000003375'
$ws.Cells.Item(38, 5).Value = '  +109.14%  '

$ws.Cells.Item(39, 2).Value = 'VeChain'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(39, 4).Value = '''0.02104'
$ws.Cells.Item(39, 5).Value = '  +1.53%  '

$ws.Cells.Item(40, 4).Value = '''7.697'
$ws.Cells.Item(40, 5).Value = '  -1.02%  '

$ws.Cells.Item(41, 4).Value = '''0.5713'
$ws.Cells.Item(41, 5).Value = '  +0.60%  '

$ws.Cells.Item(42, 2).Value = 'Algorand'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(42, 4).Value = '''0.1815'
$ws.Cells.Item(42, 5).Value = '  +2.41%  '

$ws.Cells.Item(43, 2).Value = 'Aptos'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(43, 4).Value = '''9.757'
$ws.Cells.Item(43, 5).Value = '  -0.25%  '

$ws.Cells.Item(44, 2).Value = 'EnergySwap'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(44, 4).Value = '''11.97'
$ws.Cells.Item(44, 5).Value = '  -0.15%  '

$ws.Cells.Item(45, 2).Value = 'RenderToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(45, 4).Value = '''2.210'
$ws.Cells.Item(45, 5).Value = '  -1.89%  '

$ws.Cells.Item(46, 2).Value = 'Decentraland'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(46, 4).Value = '''0.5343'
$ws.Cells.Item(46, 5).Value = '  -0.14%  '

$ws.Cells.Item(47, 2).Value = 'Cronos'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(47, 4).Value = '''0.06933'
$ws.Cells.Item(47, 5).Value = '  -1.64%  '

$ws.Cells.Item(48, 2).Value = 'MXToken'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(48, 4).Value = '''2.596'
$ws.Cells.Item(48, 5).Value = '  +2.64%  '

$ws.Cells.Item(49, 4).Value = '''1.849'
$ws.Cells.Item(49, 5).Value = '  +0.46%  '

$ws.Cells.Item(50, 2).Value = 'Quant'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(50, 4).Value = '''113.20'
$ws.Cells.Item(50, 5).Value = '  +0.77%  '

$ws.Cells.Item(51, 2).Value = 'WOONetwork'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Cells.Item(51, 4).Value = '''0.3000'
$ws.Cells.Item(51, 5).Value = '  +3.25%  '
